# Generate Report for Handback
#
# This script records the handback of the "5621a3f3-da0e-42b0-bf35-3695ad58b5bb.md"
# file for both locales (zh-cn and de-de): the status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", the target/handback
# files + handback datetime get populated, and the columns that now hold the
# longer file-name / status text are widened so the content is not clipped.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/624371383919689c2085d1993a4d33440c60560a/e2e/5621a3f3-da0e-42b0-bf35-3695ad58b5bb.md"
$targetDisplay = "5621a3f3-da0e-42b0-bf35-3695ad58b5bb.md"

# Column widths Excel snaps to a 1/6-character pixel grid, so request the
# nearest value that rounds to the desired stored width.
$wideStatusWidth = 29 + 1/6   # -> stored width ~30 (was ~17.2)
$wideFileWidth   = 39 + 1/6   # -> stored width 40 (was ~18.65 / ~21.71)

# ---------------------------------------------------------------------------
# Overview sheet: widen the per-locale status columns (E = zh-cn, F = de-de)
# and reflect the same status text (both locales share the one status string)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$overview.Columns.Item(6).ColumnWidth = $wideStatusWidth

foreach ($row in 2, 3) {
    $overview.Range("E$row").Value = "Handed back: in sync with en-US"
    $overview.Range("F$row").Value = "Handed back: in sync with en-US"
}

# ---------------------------------------------------------------------------
# Per-locale handback sheets
# ---------------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Xlf = "5621a3f3-da0e-42b0-bf35-3695ad58b5bb.d293597cea4cf024ea99325b7ae2e0558f04984c.zh-cn.xlf"; HandbackTime = "2016-08-14 17:29:24" },
    @{ Sheet = "de-de"; Xlf = "5621a3f3-da0e-42b0-bf35-3695ad58b5bb.d293597cea4cf024ea99325b7ae2e0558f04984c.de-de.xlf"; HandbackTime = "2016-08-14 17:29:34" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Widen: C = Status, I = Latest Target File, J = Latest Handback File
    $ws.Columns.Item(3).ColumnWidth = $wideStatusWidth
    $ws.Columns.Item(9).ColumnWidth = $wideFileWidth
    $ws.Columns.Item(10).ColumnWidth = $wideFileWidth

    foreach ($row in 2, 3) {
        # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Range("C$row").Value = "Handed back: in sync with en-US"

        # Latest Target File: now points at the handed-off source markdown file
        $ws.Range("I$row").Value = $targetDisplay
        $ws.Hyperlinks.Add($ws.Range("I$row"), $targetUrl, "", "", $targetDisplay) | Out-Null

        # Latest Handback File: the generated xliff that was handed back
        $ws.Range("J$row").Value = $locale.Xlf

        # Latest Handback DateTime
        $ws.Range("K$row").Value = $locale.HandbackTime
    }
}
